# Apply the "Triathlon Season" data corrections:
#  - normalize league-name casing on a few rows
#  - fix event/round text typos (Aquathon -> Aquathlon, Aqua -> Aquathlon, etc.)
#  - swap "Aquabike" -> "Standard Aquabike" where the distance was missing
#  - replace "Ironman 70.3" branding with generic "70.3"/"Ultimate" wording
#  - correct several Events/Club-name swaps that had been entered in the wrong rows
#  - add the missing club entry for Hunter League round 6 (G41)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "south coast league"
$ws.Range("A3").Value = "south coast league"
$ws.Range("A4").Value = "south coast league"
$ws.Range("F4").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("A5").Value = "south coast league"
$ws.Range("E5").Value = "Super Sprint, Sprint, Classic and Ultimate"
$ws.Range("A6").Value = "south coast league"
$ws.Range("A7").Value = "North Coast league"
$ws.Range("A8").Value = "North Coast league"
$ws.Range("F8").Value = "Super Sprint, Aquathlon, Teams"
$ws.Range("A9").Value = "North Coast league"
$ws.Range("A10").Value = "North Coast league"
$ws.Range("A11").Value = "North Coast league"
$ws.Range("F11").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("A12").Value = "North Coast league"
$ws.Range("A13").Value = "North Coast league"
$ws.Range("E14").Value = "70.3 and Sprint"
$ws.Range("F14").Value = "Aquabike"
$ws.Range("C15").Value = "Sparke Helmore Triathlon"
$ws.Range("D15").Value = "No"
$ws.Range("E15").Value = "Sprint"
$ws.Range("F15").Value = "Super Sprint, Aquabike"
$ws.Range("E17").Value = "Long Aquathlon"
$ws.Range("F17").Value = "Super Sprint Aquathlon"
$ws.Range("C18").Value = "NSW Triathlon Club Champs"
$ws.Range("D18").Value = "yes"
$ws.Range("E18").Value = "Sprint,  Standard "
$ws.Range("F18").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("E19").Value = "Super Sprint, Sprint, Classic and Ultimate"
$ws.Range("E24").Value = "70.3 and Sprint"
$ws.Range("F24").Value = "Aquabike"
$ws.Range("C25").Value = "Sparke Helmore Triathlon"
$ws.Range("D25").Value = "No"
$ws.Range("E25").Value = "Sprint"
$ws.Range("F25").Value = "Super Sprint, Aquabike"
$ws.Range("E27").Value = "Long Aquathlon"
$ws.Range("F27").Value = "Super Sprint Aquathlon"
$ws.Range("C28").Value = "NSW Triathlon Club Champs"
$ws.Range("D28").Value = "yes"
$ws.Range("E28").Value = "Sprint,  Standard "
$ws.Range("F28").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("E29").Value = "Super Sprint, Sprint, Classic and Ultimate"
$ws.Range("E30").Value = "Sprint, Olympic"
$ws.Range("E36").Value = "Sprint"
$ws.Range("F36").Value = "Super Sprint, Aquabike"
$ws.Range("G36").Value = "Forster Triathlon Club"
$ws.Range("F37").Value = "Super Sprint, Standard Aquabike"
$ws.Range("G37").Value = "Newcastle Triathlon Club"
$ws.Range("G38").Value = "Central Coast Triathlon Club"
$ws.Range("G39").Value = "Singleton Triathlon Club"
$ws.Range("G40").Value = "Maitland Triathlon Club"
$ws.Range("C41").Value = "NSW Triathlon Club Champs"
$ws.Range("E41").Value = "Sprint,  Standard "
$ws.Range("F41").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("G41").Value = "Tomaree Triathlon Club"
$ws.Range("F51").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("E52").Value = "70.3, 70.3 Aquabike, Standard, Standard Aquabike, Sprint"
$ws.Range("E54").Value = "Long Aquathlon"
$ws.Range("F55").Value = "Sprint Aquabike, Standard Aquabike, Super Sprint"
$ws.Range("E56").Value = "Super Sprint, Sprint, Classic and Ultimate"
